$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the German highlighting keys to English in the
# ticketDescriptionHighlighting column (G) for rows 2-4.

$ws.Range("G2").Value = '[
   {
      "start":143,
      "end":147,
      "key":"System"
   },
   {
      "start":104,
      "end":128,
      "key":"Error description"
   },
   {
      "start":67,
      "end":77,
      "key":"System"
   }
]'

$ws.Range("G4").Value = '[
   {
      "start":130,
      "end":165,
      "key":"Trigger"
   },
   {
      "start":37,
      "end":78,
      "key":"Error description"
   },
   {
      "start":24,
      "end":36,
      "key":"System"
   }
]'

$ws.Range("G3").Value = '[
   {
      "start":229,
      "end":297,
      "key":"Service request"
   },
   {
      "start":191,
      "end":192,
      "key":"System"
   },
   {
      "start":176,
      "end":191,
      "key":"System"
   },
   {
      "start":129,
      "end":144,
      "key":"System"
   }
]'

# Restore the original (explicit) row heights, since assigning the
# wrapped multi-line text above can trigger Excel's row auto-fit.
$ws.Rows.Item(2).RowHeight = 20
$ws.Rows.Item(3).RowHeight = 34.5
$ws.Rows.Item(4).RowHeight = 60

# Update the saved cursor/selection position to G2.
$ws.Range("G2").Select()
